$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 'Ementa atual:'
$ws.Range("C1").Value = 'Ementa modificada (dados modificados em vermelho):'

$ws.Range("B2").Value = 'LOM3111'
$ws.Range("C2").Value = 'LOM3111'

$ws.Range("A3").Value = 'Nome:'
$ws.Range("B3").Value = ' Processamento de Cerâmicas Experimental'
$ws.Range("C3").Value = ' Processamento de Cerâmicas Experimental'

$ws.Range("A4").Value = 'Name:'
$ws.Range("B4").Value = 'Experimental Ceramics Processing'
$ws.Range("C4").Value = 'Experimental Ceramics Processing'

$ws.Range("A5").Value = 'Créditos-aula:'
$ws.Range("B5").Value = '4'
$ws.Range("C5").Value = '4'

$ws.Range("A6").Value = 'Créditos-trabalho'
$ws.Range("B6").Value = '0'
$ws.Range("C6").Value = '0'

$ws.Range("A7").Value = 'Carga horária:'
$ws.Range("B7").Value = '60 h'
$ws.Range("C7").Value = '60 h'

$ws.Range("A8").Value = 'Ativação:'
$ws.Range("B8").Value = '01/01/2022'
$ws.Range("C8").Value = '01/01/2022'

$ws.Range("A9").Value = 'Semestre ideal:'
$ws.Range("B9").Value = 'EM-7'
$ws.Range("C9").Value = 'EM-7'

$ws.Range("A10").Value = 'Objetivos:'
$ws.Range("B10").Value = 'Esta disciplina faz parte da formação do engenheiro de materiais, contribuindo para gerar competências gerais e específicas.Incentivar trabalhos em grupo, com ênfase na visão integrada sobre os aspectos abordados na disciplina.Desenvolver habilidades práticas, bem como conhecer e operar equipamentosPromover a comunicação nas formas escrita, oral e gráfica, além de trabalhos em grupos.Relacionar esta disciplina com outras da grade do curso, tanto com as de formação específica quanto às de formação geral.'
$ws.Range("C10").Value = 'Esta disciplina faz parte da formação do engenheiro de materiais, contribuindo para gerar competências gerais e específicas.Incentivar trabalhos em grupo, com ênfase na visão integrada sobre os aspectos abordados na disciplina.Desenvolver habilidades práticas, bem como conhecer e operar equipamentosPromover a comunicação nas formas escrita, oral e gráfica, além de trabalhos em grupos.Relacionar esta disciplina com outras da grade do curso, tanto com as de formação específica quanto às de formação geral.'

$ws.Range("A11").Value = 'Objectives:'

$ws.Range("A12").Value = 'Docentes responsáveis:'

$ws.Range("B13").Value = '5983729 - Fernando Vernilli Junior'
$ws.Range("C13").Value = '5983729 - Fernando Vernilli Junior'

$ws.Range("B14").Value = '1922320 - Sebastiao Ribeiro'
$ws.Range("C14").Value = '1922320 - Sebastiao Ribeiro'

$ws.Range("A15").Value = 'Programa resumido:'
$ws.Range("B15").Value = '1 – Introdução, 2 – Moagem e Análise granulométrica, 3 - Conformação por via seca, 4 - Secagem e queima (sinterização), 5 – Determinação da massa específica aparente, absorção de água e porosidade,  6 - Ensaios de resistência mecânica. 7 - Preparação e caracterização de suspensões cerâmicas, 8 - Confecção de moldes de gesso para conformação por via liquida, 9 - Conformação por via líquidda (colagem de barbotine), 10 - Limite de plasticidade e liquidez, 11 - Preparação e conformação de massas cerâmicas pastosas, 12 - Preparo de superfícies cerâmicas.'
$ws.Range("C15").Value = '1 – Introdução, 2 – Moagem e Análise granulométrica, 3 - Conformação por via seca, 4 - Secagem e queima (sinterização), 5 – Determinação da massa específica aparente, absorção de água e porosidade,  6 - Ensaios de resistência mecânica. 7 - Preparação e caracterização de suspensões cerâmicas, 8 - Confecção de moldes de gesso para conformação por via liquida, 9 - Conformação por via líquidda (colagem de barbotine), 10 - Limite de plasticidade e liquidez, 11 - Preparação e conformação de massas cerâmicas pastosas, 12 - Preparo de superfícies cerâmicas.'

$ws.Range("A16").Value = 'Short syllabus:'

$ws.Range("A17").Value = 'Programa:'
$ws.Range("B17").Value = '1 Introdução, 2 – Moagem e Análise granulométrica – parâmetros de moagem em moinhos de bola, peneiração, 3 - Conformação por via seca: prensagem uniaxial, construção de curvas de compactação, 4 - Secagem e queima (sinterização): elaboração de curvas de secagem, avaliação da retração e densificação, 5 – Determinação da massa específica aparente, absorção de água e porosidade, 6 - Ensaios de resistência mecânica à flexão e compressão diametral – modulo de Weibull,  7 - Preparação e caracterização de suspensões cerâmicas:  medida de massa especifica, viscosidade e estabilidade,  8 - Confecção de moldes de gesso para conformação por via liquida (colagem de barbotina), 9 - Conformação por colagem de barbotina, 10 - Limite de plasticidade e liquidez, 11 - Preparação e conformação de massas cerâmicas  pastosas: extrusão convencional (maromba) e conformação aditiva (impressão 3D), 12 - Preparo de superfícies cerâmicas:  Vidragem, retificação, lixamento e polimento'
$ws.Range("C17").Value = '1 Introdução, 2 – Moagem e Análise granulométrica – parâmetros de moagem em moinhos de bola, peneiração, 3 - Conformação por via seca: prensagem uniaxial, construção de curvas de compactação, 4 - Secagem e queima (sinterização): elaboração de curvas de secagem, avaliação da retração e densificação, 5 – Determinação da massa específica aparente, absorção de água e porosidade, 6 - Ensaios de resistência mecânica à flexão e compressão diametral – modulo de Weibull,  7 - Preparação e caracterização de suspensões cerâmicas:  medida de massa especifica, viscosidade e estabilidade,  8 - Confecção de moldes de gesso para conformação por via liquida (colagem de barbotina), 9 - Conformação por colagem de barbotina, 10 - Limite de plasticidade e liquidez, 11 - Preparação e conformação de massas cerâmicas  pastosas: extrusão convencional (maromba) e conformação aditiva (impressão 3D), 12 - Preparo de superfícies cerâmicas:  Vidragem, retificação, lixamento e polimento'

$ws.Range("A18").Value = 'Syllabus:'

$ws.Range("A19").Value = 'Avaliação:'

$ws.Range("A20").Value = 'Método:'
$ws.Range("B20").Value = 'Avaliação individual do comportamento do aluno frente aos trabalhos experimentais (AC), Relatórios sobre os testes experimentais (MAR) e prova experimental final (PE).'
$ws.Range("C20").Value = 'Avaliação individual do comportamento do aluno frente aos trabalhos experimentais (AC), Relatórios sobre os testes experimentais (MAR) e prova experimental final (PE).'

$ws.Range("A21").Value = 'Critério:'
$ws.Range("B21").Value = '1 - Média aritmética das notas dos relatórios, com peso 1 (MAR), 2 – avaliação comportamental, peso 1 (AC) e uma prova experimental (PE), no final do semestre letivo, com peso 2.A nota final (NF) será calculada pela equação (MAR+AC+2PE)/4. NF igual ou superior a 5: aprovação.'
$ws.Range("C21").Value = '1 - Média aritmética das notas dos relatórios, com peso 1 (MAR), 2 – avaliação comportamental, peso 1 (AC) e uma prova experimental (PE), no final do semestre letivo, com peso 2.A nota final (NF) será calculada pela equação (MAR+AC+2PE)/4. NF igual ou superior a 5: aprovação.'

$ws.Range("A22").Value = 'Norma de recuperação:'
$ws.Range("B22").Value = 'Devido à característica da disciplina não'
$ws.Range("C22").Value = 'Devido à característica da disciplina não'

$ws.Range("A23").Value = 'Bibliografia:'
$ws.Range("B23").Value = '1. Dispersão e empacotamento de partículas, Fazendo Arte Editorial. Ivone R de Oliveira e co-autores, 2000, 2. Norton, F.H. Introdução à Tecnologia Cerâmica, Ed. Edgard Blucher, 1973, 3. Reed, J.S. Principles of Ceramics Processing, John Wiley, 1988, 4. Rahaman, M. N. Ceramic Processing and Sintering. 2st Editon, 2003, 5.Van Vlack, L.M. Propriedades dos Materiais Cerâmicos, Ed. Edgard Blucher, 1973, 6. Ceramic Materials: Science and Engineering, C. Barry Carter, M. Grant Norton 2nd ed., 2013, 7. Fundamentals of Ceramic Powder Processing and Synthesis: Terry A. Ring,  8. R. A. Terpsta, P. P. A. C. Pex, A. H. de Vries, Ceramic Processing, Edited: R. A. Terpsta, P. P. A. C. Pex, A. H. de Vries, 1995. 9. M. F. Ashby, D.R. H. Jones, Engenharia de Materiais, Volume II, 3° edição, Elsevier, p.436, 2007,10 – Normas ASTM, ISSO e ABNT11. Artigos da literatura especializada,'
$ws.Range("C23").Value = '1. Dispersão e empacotamento de partículas, Fazendo Arte Editorial. Ivone R de Oliveira e co-autores, 2000, 2. Norton, F.H. Introdução à Tecnologia Cerâmica, Ed. Edgard Blucher, 1973, 3. Reed, J.S. Principles of Ceramics Processing, John Wiley, 1988, 4. Rahaman, M. N. Ceramic Processing and Sintering. 2st Editon, 2003, 5.Van Vlack, L.M. Propriedades dos Materiais Cerâmicos, Ed. Edgard Blucher, 1973, 6. Ceramic Materials: Science and Engineering, C. Barry Carter, M. Grant Norton 2nd ed., 2013, 7. Fundamentals of Ceramic Powder Processing and Synthesis: Terry A. Ring,  8. R. A. Terpsta, P. P. A. C. Pex, A. H. de Vries, Ceramic Processing, Edited: R. A. Terpsta, P. P. A. C. Pex, A. H. de Vries, 1995. 9. M. F. Ashby, D.R. H. Jones, Engenharia de Materiais, Volume II, 3° edição, Elsevier, p.436, 2007,10 – Normas ASTM, ISSO e ABNT11. Artigos da literatura especializada,'

$ws.Range("A24").Value = 'Requisitos:'

$ws.Range("B25").Value = 'LOM3073 -  Processamento de Cerâmicas  (Requisito fraco)
'
$ws.Range("C25").Value = 'LOM3073 -  Processamento de Cerâmicas  (Requisito fraco)
'

$ws.Range("A13").ClearContents()
$ws.Range("A14").ClearContents()
$ws.Range("B18").ClearContents()
$ws.Range("C18").ClearContents()
$ws.Range("B19").ClearContents()
$ws.Range("C19").ClearContents()

$ws.Rows.Item(10).RowHeight = 60
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 120
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(23).RowHeight = 120
$ws.Rows.Item(25).RowHeight = 30

$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).AutoFit()
$ws.Rows.Item(6).AutoFit()
$ws.Rows.Item(7).AutoFit()
$ws.Rows.Item(8).AutoFit()
$ws.Rows.Item(9).AutoFit()
$ws.Rows.Item(12).AutoFit()
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).AutoFit()
$ws.Rows.Item(19).AutoFit()
$ws.Rows.Item(24).AutoFit()
